# Auto-generated Word COM-interop script to update division problems
$d = $word.ActiveDocument

$d.Content.Find.Execute("600÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "802÷2=", 2) | Out-Null
$d.Content.Find.Execute("296÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "127÷8=", 2) | Out-Null
$d.Content.Find.Execute("843÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "964÷3=", 2) | Out-Null
$d.Content.Find.Execute("334÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "748÷4=", 2) | Out-Null
$d.Content.Find.Execute("648÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "355÷9=", 2) | Out-Null
$d.Content.Find.Execute("423÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "281÷7=", 2) | Out-Null
$d.Content.Find.Execute("747÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "370÷8=", 2) | Out-Null
$d.Content.Find.Execute("277÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "851÷5=", 2) | Out-Null
$d.Content.Find.Execute("331÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "640÷7=", 2) | Out-Null
$d.Content.Find.Execute("565÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "435÷8=", 2) | Out-Null
$d.Content.Find.Execute("223÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "210÷9=", 2) | Out-Null
$d.Content.Find.Execute("597÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "298÷5=", 2) | Out-Null
$d.Content.Find.Execute("254÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "247÷6=", 2) | Out-Null
$d.Content.Find.Execute("944÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "535÷7=", 2) | Out-Null
$d.Content.Find.Execute("585÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "129÷7=", 2) | Out-Null
$d.Content.Find.Execute("421÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "217÷2=", 2) | Out-Null
$d.Content.Find.Execute("151÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "942÷7=", 2) | Out-Null
$d.Content.Find.Execute("492÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "578÷8=", 2) | Out-Null
$d.Content.Find.Execute("492÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "555÷9=", 2) | Out-Null
$d.Content.Find.Execute("185÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "828÷6=", 2) | Out-Null
$d.Content.Find.Execute("969÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "711÷9=", 2) | Out-Null
$d.Content.Find.Execute("796÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "244÷9=", 2) | Out-Null
$d.Content.Find.Execute("481÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "303÷6=", 2) | Out-Null
$d.Content.Find.Execute("332÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "768÷6=", 2) | Out-Null
$d.Content.Find.Execute("684÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "969÷2=", 2) | Out-Null

Write-Host "Replacement complete"
